# Refresh the quoted crypto symbol list (Price / Volume(1h) columns) to
# the latest scraped snapshot. Values are stored as literal text (the
# source sheet keeps numbers/percentages as text, e.g. "307.93",
# "-1.20%"), so a plain `.Value = "307.93"` assignment would be wrongly
# re-interpreted by Excel as a number/percentage. Instead we assign a
# text-producing formula (`="307.93"`) and then Copy / PasteSpecial
# (values-only) that same cell back onto itself -- this commits the
# formula's cached text result as a plain literal without Excel's normal
# type-sniffing on direct value assignment, and without touching the
# cell's number format/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="307.93"'
$ws.Range("E2").Formula = '="-1.20%"'
$rng = $ws.Range("D2:E2")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D3").Formula = '="40.83"'
$ws.Range("E3").Formula = '="-0.55%"'
$rng = $ws.Range("D3:E3")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D4").Formula = '="5.041"'
$ws.Range("E4").Formula = '="-1.32%"'
$rng = $ws.Range("D4:E4")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D5").Formula = '="0.07635"'
$ws.Range("E5").Formula = '="-2.96%"'
$rng = $ws.Range("D5:E5")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D6").Formula = '="4.263"'
$ws.Range("E6").Formula = '="-1.79%"'
$rng = $ws.Range("D6:E6")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D7").Formula = '="1.619"'
$ws.Range("E7").Formula = '="-4.10%"'
$rng = $ws.Range("D7:E7")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D9").Formula = '="0.9095"'
$ws.Range("E9").Formula = '="-1.51%"'
$rng = $ws.Range("D9:E9")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("E10").Formula = '="-7.52%"'
$rng = $ws.Range("E10")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("E11").Formula = '="-0.58%"'
$rng = $ws.Range("E11")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D12").Formula = '="0.09058"'
$ws.Range("E12").Formula = '="-0.15%"'
$rng = $ws.Range("D12:E12")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D13").Formula = '="0.04315"'
$ws.Range("E13").Formula = '="-1.85%"'
$rng = $ws.Range("D13:E13")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D14").Formula = '="0.1053"'
$ws.Range("E14").Formula = '="-0.49%"'
$rng = $ws.Range("D14:E14")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D15").Formula = '="0.001252"'
$ws.Range("E15").Formula = '="-0.69%"'
$rng = $ws.Range("D15:E15")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D16").Formula = '="0.005797"'
$ws.Range("E16").Formula = '="-2.54%"'
$rng = $ws.Range("D16:E16")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("E17").Formula = '="-0.67%"'
$rng = $ws.Range("E17")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("E18").Formula = '="-2.95%"'
$rng = $ws.Range("E18")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D19").Formula = '="6.795"'
$ws.Range("E19").Formula = '="-5.46%"'
$rng = $ws.Range("D19:E19")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D20").Formula = '="0.1357"'
$ws.Range("E20").Formula = '="-1.30%"'
$rng = $ws.Range("D20:E20")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D21").Formula = '="0.2723"'
$ws.Range("E21").Formula = '="-2.76%"'
$rng = $ws.Range("D21:E21")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D22").Formula = '="0.04156"'
$ws.Range("E22").Formula = '="-0.12%"'
$rng = $ws.Range("D22:E22")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("E23").Formula = '="0.27%"'
$rng = $ws.Range("E23")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D24").Formula = '="0.004089"'
$ws.Range("E24").Formula = '="-1.55%"'
$rng = $ws.Range("D24:E24")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D25").Formula = '="0.0001300"'
$ws.Range("E25").Formula = '="6.09%"'
$rng = $ws.Range("D25:E25")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D26").Formula = '="0.0003007"'
$ws.Range("E26").Formula = '="0.56%"'
$rng = $ws.Range("D26:E26")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D38").Formula = '="0.02403"'
$ws.Range("E38").Formula = '="-2.17%"'
$rng = $ws.Range("D38:E38")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D39").Formula = '="0.05165"'
$ws.Range("E39").Formula = '="-2.82%"'
$rng = $ws.Range("D39:E39")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D40").Formula = '="0.007777"'
$ws.Range("E40").Formula = '="-3.14%"'
$rng = $ws.Range("D40:E40")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D41").Formula = '="0.1308"'
$ws.Range("E41").Formula = '="-3.48%"'
$rng = $ws.Range("D41:E41")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D42").Formula = '="0.007076"'
$ws.Range("E42").Formula = '="-6.11%"'
$rng = $ws.Range("D42:E42")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D43").Formula = '="0.001917"'
$ws.Range("E43").Formula = '="-3.64%"'
$rng = $ws.Range("D43:E43")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D44").Formula = '="0.007470"'
$ws.Range("E44").Formula = '="-8.82%"'
$rng = $ws.Range("D44:E44")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D45").Formula = '="0.3345"'
$ws.Range("E45").Formula = '="7.73%"'
$rng = $ws.Range("D45:E45")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D46").Formula = '="0.00006348"'
$ws.Range("E46").Formula = '="-6.18%"'
$rng = $ws.Range("D46:E46")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D47").Formula = '="0.00000000750"'
$ws.Range("E47").Formula = '="-0.39%"'
$rng = $ws.Range("D47:E47")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D48").Formula = '="0.004400"'
$ws.Range("E48").Formula = '="6.90%"'
$rng = $ws.Range("D48:E48")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D49").Formula = '="0.006211"'
$ws.Range("E49").Formula = '="81.87%"'
$rng = $ws.Range("D49:E49")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D50").Formula = '="0.00002100"'
$ws.Range("E50").Formula = '="-0.39%"'
$rng = $ws.Range("D50:E50")
$rng.Copy()
$rng.PasteSpecial(-4163)

$ws.Range("D51").Formula = '="0.0002000"'
$ws.Range("E51").Formula = '="-0.39%"'
$rng = $ws.Range("D51:E51")
$rng.Copy()
$rng.PasteSpecial(-4163)

$excel.CutCopyMode = $false
